# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) data block in rows
# 16-50 gets reordered so the most recent period (2009) is now first and
# the oldest (1711) is now last - i.e. the 35-row block is reversed in
# place, row by row, while every other column (B, C, D, G, H, I, J) and
# all cell styles stay exactly where they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow = 50
$count = $lastRow - $firstRow + 1

# Column numbers: E = 5 (Periodo Mora), F = 6 (Valor Mora)
$periodCol = 5
$valueCol = 6

# 1) Snapshot the current values for the block
$periods = @()
$values = @()
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $periods += $ws.Cells.Item($r, $periodCol).Value()
    $values += $ws.Cells.Item($r, $valueCol).Value()
}

# 2) Write the snapshot back in reverse order
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $srcIndex = $count - 1 - $i
    $ws.Cells.Item($r, $periodCol).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, $valueCol).Value = $values[$srcIndex]
}
